# Add a new "Sheet2" after "Sheet1", populate it with the twizzlers survey
# table, size its columns, and leave the final on-screen selection matching
# the authored file (Sheet1: A1:D5 selected, Sheet2: B6 selected/active).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New worksheet, inserted immediately after Sheet1 (becomes Sheet2 / rId2).
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Header row.
$ws2.Range("A1").Value = " "
$ws2.Range("B1").Value = "How many twizzlers did you eat?"
$ws2.Range("C1").Value = "Correct"
$ws2.Range("D1").Value = "Comment"

# Row 2.
$ws2.Range("A2").Value = "A"
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = "N"
$ws2.Range("D2").Value = "Obviously it's not red."

# Row 3.
$ws2.Range("A3").Value = "B"
$ws2.Range("B3").Value = 3
$ws2.Range("C3").Value = "Y"
$ws2.Range("D3").Value = "You're brilliant!"

# Row 4.
$ws2.Range("A4").Value = "C"
$ws2.Range("B4").Value = 4
$ws2.Range("C4").Value = "N"
$ws2.Range("D4").Value = "You're stupid."

# Row 5.
$ws2.Range("A5").Value = "D"
$ws2.Range("B5").Value = 5
$ws2.Range("C5").Value = "N"
$ws2.Range("D5").Value = "So close. But not really that close."

# Column widths (character units; converted internally to the stored
# spreadsheetML "width" - closest achievable values to 6.71 / 35 / 45.29).
$ws2.Columns.Item(1).ColumnWidth = 5.833333333333333
$ws2.Columns.Item(2).ColumnWidth = 34.166666666666664
$ws2.Columns.Item(4).ColumnWidth = 44.5

# Sheet1 keeps its data untouched but its selection becomes the full table
# (A1:D5) and it is no longer the active tab.
$ws1.Range("A1:D5").Select()

# Sheet2 ends up the active sheet/tab, with B6 selected.
$ws2.Range("B6").Select()
